# Applies the "PO Forecast" update:
#  1. Rename "Requested quantity" header -> "Weekly_PO_Qty"  (Weekly Quantity sheet, B1)
#  2. Rename "Requested quantity" header -> "Monthly_PO_Qty" (Monthly Trend sheet, B1)
#  3. Add a new "PO Forecast" worksheet (as the 3rd / last sheet) with
#     forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity header rename ---------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend header rename -----------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. New "PO Forecast" sheet ---------------------------------------------
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"
# Move it so it becomes the last sheet (after the current last sheet).
$wsForecast.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
# Re-resolve the worksheet handle by name (post-move, safest reference).
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
  @(45340.99999999999, 0, -24.99455397570485, 10.62071325473459),
  @(45375.99999999999, 0, -20.47983989709691, 18.65580657540976),
  @(45389.99999999999, 1, -15.6039786603724, 18.44912068674154),
  @(45403.99999999999, 4, -15.21015836924235, 21.49276532389167),
  @(45431.99999999999, 8, -12.32949662048703, 26.28933099593345),
  @(45445.99999999999, 11, -8.667259161093442, 29.76387306168566),
  @(45452.99999999999, 12, -7.179065210539243, 29.9387961614021),
  @(45459.99999999999, 13, -5.898033273167417, 29.4288878779542),
  @(45487.99999999999, 18, -0.6866365700226914, 35.96317941957778),
  @(45508.99999999999, 21, 2.784222621341873, 39.18685539886252),
  @(45515.99999999999, 22, 3.961507462269044, 40.20472481601283),
  @(45529.99999999999, 25, 6.153853271883044, 43.44413564362883),
  @(45550.99999999999, 28, 10.66769351322982, 46.71104809728499),
  @(45557.99999999999, 29, 10.44349410680329, 47.64813436945327),
  @(45578.99999999999, 33, 15.24794846137248, 50.6554933053426),
  @(45585.99999999999, 34, 14.61916861495184, 51.80723644466568),
  @(45592.99999999999, 35, 18.55898793704416, 54.35337078168578),
  @(45599.99999999999, 36, 17.37088270864909, 53.81962592766918),
  @(45606.99999999999, 37, 19.99887717656677, 54.93117328522171),
  @(45613.99999999999, 38, 19.34881595304724, 56.56173702966356),
  @(45620.99999999999, 40, 20.25873045291897, 58.75124351400969),
  @(45627.99999999999, 41, 21.58176011069866, 58.19818105010038),
  @(45634.99999999999, 42, 24.34557266132923, 59.77442175118065),
  @(45641.99999999999, 43, 24.03857455571851, 61.75680816331612)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}
$lastRow = $row - 1

# --- Formatting: reuse the existing header / date-column styles from the
# "Weekly Quantity" sheet (same bold+border+centered header, same
# yyyy-mm-dd date format on column A) so the new sheet matches the look of
# the other two sheets exactly. Copy + PasteSpecial(formats only) so cell
# VALUES already written above are left untouched.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A$lastRow").PasteSpecial(-4122)

$wsForecast.Application.CutCopyMode = $false

Write-Output "PO Forecast sheet created with $($data.Count) data rows."
